$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# New identifiers used by this handoff report regeneration.
# -----------------------------------------------------------------
$newUuid = "7c52a147-6b73-4930-9b08-b0b049f4e3e2"
$newUuid2 = "ffff8499a01e-c57e-47fe-8bdd-b6aa98e8fe4b"
$newHash = "698c69b311bc70b073a0eece03175770633e26df"

$newMdName = "$newUuid.md"
$newMdName2 = "$newUuid2.md"
$newZhXlf = "$newUuid.$newHash.zh-cn.xlf"
$newDeXlf = "$newUuid.$newHash.de-de.xlf"

$newHandoffDate = "2016-35-14 09:35:50"
$newZhDatetime = "2016-03-14 09:35:42"
$newDeDatetime = "2016-03-14 09:35:50"

# =========================================================================
# Sheet "Overview"
# =========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

# -- Row 2: refresh the existing handoff entry with the new file name/date --
$wsOverview.Range("A2").Hyperlinks.Delete()
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = $newHandoffDate
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/aad113a7ea495081e034b20f1230e52c3286f1d6/e2e/$newMdName",
    "",
    "",
    $newMdName
) | Out-Null

# -- Row 3: new handoff entry appended by the report generator --
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = $newHandoffDate
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/aad113a7ea495081e034b20f1230e52c3286f1d6/e2e/$newMdName2",
    "",
    "",
    $newMdName2
) | Out-Null

# =========================================================================
# Sheet "zh-cn"
# =========================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

# -- Row 2: refresh with the new source/handoff file + timestamp --
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("B2").Hyperlinks.Delete()
$wsZh.Range("D2").Hyperlinks.Delete()

$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("E2").Value = $newZhDatetime
$wsZh.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H2").Value = "0001-01-01 00:00:00"
$wsZh.Range("I2").Value = "Include"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/aad113a7ea495081e034b20f1230e52c3286f1d6/e2e/$newMdName",
    "",
    "",
    $newMdName
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("B2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/aad113a7ea495081e034b20f1230e52c3286f1d6/e2e/$newMdName",
    "",
    "",
    ".md"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f6d844cb9f82f5193a8e60933d0ffbcca6f04482/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/$newZhXlf",
    "",
    "",
    $newZhXlf
) | Out-Null

# -- Row 3: new entry appended by the report generator --
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("E3").Value = $newZhDatetime
$wsZh.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("I3").Value = "Include"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/aad113a7ea495081e034b20f1230e52c3286f1d6/e2e/$newMdName2",
    "",
    "",
    $newMdName2
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("B3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/aad113a7ea495081e034b20f1230e52c3286f1d6/e2e/$newMdName2",
    "",
    "",
    ".md"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f6d844cb9f82f5193a8e60933d0ffbcca6f04482/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/$newZhXlf",
    "",
    "",
    $newZhXlf
) | Out-Null

# =========================================================================
# Sheet "de-de"
# =========================================================================
$wsDe = $wb.Worksheets.Item("de-de")

# -- Row 2: refresh with the new source/handoff file + timestamp --
$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("B2").Hyperlinks.Delete()
$wsDe.Range("D2").Hyperlinks.Delete()

$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("E2").Value = $newDeDatetime
$wsDe.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H2").Value = "0001-01-01 00:00:00"
$wsDe.Range("I2").Value = "Include"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/aad113a7ea495081e034b20f1230e52c3286f1d6/e2e/$newMdName",
    "",
    "",
    $newMdName
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("B2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/aad113a7ea495081e034b20f1230e52c3286f1d6/e2e/$newMdName",
    "",
    "",
    ".md"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/30214c4ab0c6d2c71fbafa1fc8b4114eae79e88c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/$newDeXlf",
    "",
    "",
    $newDeXlf
) | Out-Null

# -- Row 3: new entry appended by the report generator --
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("E3").Value = $newDeDatetime
$wsDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("I3").Value = "Include"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/aad113a7ea495081e034b20f1230e52c3286f1d6/e2e/$newMdName2",
    "",
    "",
    $newMdName2
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("B3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/aad113a7ea495081e034b20f1230e52c3286f1d6/e2e/$newMdName2",
    "",
    "",
    ".md"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/30214c4ab0c6d2c71fbafa1fc8b4114eae79e88c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/$newDeXlf",
    "",
    "",
    $newDeXlf
) | Out-Null
